# Auto-generated edit script applying numeric corrections to the
# Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 11170.143
$ws.Cells.Item(43, 9).Value = 24519.8
$ws.Cells.Item(43, 10).Value = 3753.6667
$ws.Cells.Item(43, 11).Value = 24519.8
$ws.Cells.Item(43, 12).Value = 3753.6667
$ws.Cells.Item(43, 13).Value = -24450.8
$ws.Cells.Item(43, 14).Value = -3891.6667

$ws.Cells.Item(51, 8).Value = 2742.8696
$ws.Cells.Item(51, 10).Value = 3617.2
$ws.Cells.Item(51, 12).Value = 3617.2
$ws.Cells.Item(51, 14).Value = -4585.2

$ws.Cells.Item(111, 8).Value = 8173.25
$ws.Cells.Item(111, 9).Value = 13324.25
$ws.Cells.Item(111, 10).Value = 3022.25
$ws.Cells.Item(111, 11).Value = 39972.75
$ws.Cells.Item(111, 12).Value = 9066.75
$ws.Cells.Item(111, 13).Value = -36905.75
$ws.Cells.Item(111, 14).Value = -15200.75

$ws.Cells.Item(132, 8).Value = 1743.0385
$ws.Cells.Item(132, 9).Value = 1479.4584
$ws.Cells.Item(132, 11).Value = 4438.3752
$ws.Cells.Item(132, 13).Value = -1908.3752

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 12572.5
$ws.Cells.Item(110, 9).Value = 14336.071
$ws.Cells.Item(110, 11).Value = 14336.071
$ws.Cells.Item(110, 13).Value = -12291.071

$ws.Cells.Item(132, 8).Value = 2498.2856
$ws.Cells.Item(132, 9).Value = 2237.3215
$ws.Cells.Item(132, 10).Value = 3542.1428
$ws.Cells.Item(132, 11).Value = 6711.9645
$ws.Cells.Item(132, 12).Value = 10626.4284
$ws.Cells.Item(132, 13).Value = -4181.9645
$ws.Cells.Item(132, 14).Value = -15686.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3364.4424
$ws.Cells.Item(20, 9).Value = 2843.2903
$ws.Cells.Item(20, 10).Value = 4133.7617
$ws.Cells.Item(20, 11).Value = 2843.2903
$ws.Cells.Item(20, 12).Value = 4133.7617
$ws.Cells.Item(20, 13).Value = -2596.2903
$ws.Cells.Item(20, 14).Value = -4627.7617

$ws.Cells.Item(94, 8).Value = 3022.5833
$ws.Cells.Item(94, 9).Value = 1726.1
$ws.Cells.Item(94, 11).Value = 1726.1
$ws.Cells.Item(94, 13).Value = -1275.1

$ws.Cells.Item(134, 8).Value = 1924.2817
$ws.Cells.Item(134, 9).Value = 1907.7059
$ws.Cells.Item(134, 11).Value = 5723.1177
$ws.Cells.Item(134, 13).Value = -3188.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 70178.13
$ws.Cells.Item(31, 9).Value = 93859.91
$ws.Cells.Item(31, 11).Value = 93859.91
$ws.Cells.Item(31, 13).Value = -93564.91

$ws.Cells.Item(34, 8).Value = 70178.13
$ws.Cells.Item(34, 9).Value = 93859.91
$ws.Cells.Item(34, 11).Value = 93859.91
$ws.Cells.Item(34, 13).Value = -93657.91

$ws.Cells.Item(58, 8).Value = 3604
$ws.Cells.Item(58, 9).Value = 3276.2
$ws.Cells.Item(58, 10).Value = 4150.3335
$ws.Cells.Item(58, 11).Value = 3276.2
$ws.Cells.Item(58, 12).Value = 4150.3335
$ws.Cells.Item(58, 13).Value = -3073.2
$ws.Cells.Item(58, 14).Value = -4556.3335

$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = $null
$ws.Cells.Item(86, 14).Value = $null

$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = $null
$ws.Cells.Item(89, 14).Value = $null

$ws.Cells.Item(132, 8).Value = 3395.6667
$ws.Cells.Item(132, 9).Value = 3652.923
$ws.Cells.Item(132, 11).Value = 10958.769
$ws.Cells.Item(132, 13).Value = -8428.769

$ws.Cells.Item(136, 8).Value = 3604
$ws.Cells.Item(136, 9).Value = 3276.2
$ws.Cells.Item(136, 10).Value = 4150.3335
$ws.Cells.Item(136, 11).Value = 9828.599999999999
$ws.Cells.Item(136, 12).Value = 12451.0005
$ws.Cells.Item(136, 13).Value = -7278.599999999999
$ws.Cells.Item(136, 14).Value = -17551.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 97.416664
$ws.Cells.Item(2, 9).Value = 104.666664
$ws.Cells.Item(2, 10).Value = 75.666664
$ws.Cells.Item(2, 11).Value = 627.999984
$ws.Cells.Item(2, 12).Value = 453.999984
$ws.Cells.Item(2, 13).Value = -514.999984
$ws.Cells.Item(2, 14).Value = -679.999984

$ws.Cells.Item(7, 8).Value = 337.5
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 14).Value = $null

$ws.Cells.Item(32, 8).Value = 3547.625
$ws.Cells.Item(32, 10).Value = 3411.5715
$ws.Cells.Item(32, 12).Value = 10234.7145
$ws.Cells.Item(32, 14).Value = -10800.7145

$ws.Cells.Item(121, 8).Value = 2738
$ws.Cells.Item(121, 10).Value = 915
$ws.Cells.Item(121, 12).Value = 2745
$ws.Cells.Item(121, 14).Value = -5365

$ws.Cells.Item(131, 8).Value = 23707.479
$ws.Cells.Item(131, 9).Value = 112875.22
$ws.Cells.Item(131, 10).Value = 2018.027
$ws.Cells.Item(131, 11).Value = 338625.66
$ws.Cells.Item(131, 12).Value = 6054.081
$ws.Cells.Item(131, 13).Value = -333585.66
$ws.Cells.Item(131, 14).Value = -16134.081

$ws.Cells.Item(137, 8).Value = 3249.2144
$ws.Cells.Item(137, 10).Value = 3939.1428
$ws.Cells.Item(137, 12).Value = 11817.4284
$ws.Cells.Item(137, 14).Value = -22017.4284

$ws.Cells.Item(139, 8).Value = 2108.8125
$ws.Cells.Item(139, 9).Value = 1811.8334
$ws.Cells.Item(139, 10).Value = 2999.75
$ws.Cells.Item(139, 11).Value = 5435.5002
$ws.Cells.Item(139, 12).Value = 8999.25
$ws.Cells.Item(139, 13).Value = -295.5002000000004
$ws.Cells.Item(139, 14).Value = -19279.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(4, 8).Value = 5000
$ws.Cells.Item(4, 10).Value = 5000
$ws.Cells.Item(4, 12).Value = 5000
$ws.Cells.Item(4, 14).Value = -5224

$ws.Cells.Item(70, 8).Value = 3806.8823
$ws.Cells.Item(70, 9).Value = 3721.7144
$ws.Cells.Item(70, 11).Value = 3721.7144
$ws.Cells.Item(70, 13).Value = -3451.7144

$ws.Cells.Item(73, 8).Value = 3806.8823
$ws.Cells.Item(73, 9).Value = 3721.7144
$ws.Cells.Item(73, 11).Value = 3721.7144
$ws.Cells.Item(73, 13).Value = -2785.7144

$ws.Cells.Item(122, 8).Value = 2280.3257
$ws.Cells.Item(122, 9).Value = 2065.158
$ws.Cells.Item(122, 11).Value = 6195.474
$ws.Cells.Item(122, 13).Value = -3745.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4247.1665
$ws.Cells.Item(61, 9).Value = 4085.5334
$ws.Cells.Item(61, 10).Value = 5055.3335
$ws.Cells.Item(61, 11).Value = 4085.5334
$ws.Cells.Item(61, 12).Value = 5055.3335
$ws.Cells.Item(61, 13).Value = -3883.5334
$ws.Cells.Item(61, 14).Value = -5459.3335

$ws.Cells.Item(113, 8).Value = 4247.1665
$ws.Cells.Item(113, 9).Value = 4085.5334
$ws.Cells.Item(113, 10).Value = 5055.3335
$ws.Cells.Item(113, 11).Value = 4085.5334
$ws.Cells.Item(113, 12).Value = 5055.3335
$ws.Cells.Item(113, 13).Value = -1915.5334
$ws.Cells.Item(113, 14).Value = -9395.333500000001

$ws.Cells.Item(122, 8).Value = 3911.8
$ws.Cells.Item(122, 9).Value = 3329.2354
$ws.Cells.Item(122, 11).Value = 9987.706200000001
$ws.Cells.Item(122, 13).Value = -7537.706200000001

$ws.Cells.Item(140, 8).Value = 162969.75
$ws.Cells.Item(140, 10).Value = 194000
$ws.Cells.Item(140, 12).Value = 194000
$ws.Cells.Item(140, 14).Value = -204360
